$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# A new case (Next.js front-end job) was scraped at 2025-09-06 12:31:29 and is
# now the top row; every previously-listed row shifts down by one and gets
# its "fetched at" timestamp refreshed to match this run.
$ws.Rows.Item(2).Insert()

$ws.Range("A2").Value = "2025-09-06 12:31:29"
$ws.Range("B2").Value = "初回 【急募】フロントエンド Next.js システム開発"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5388329"
$ws.Range("G2").Value = 243
$ws.Range("H2").Value = "🔥Next.js ◆開発,システム開発"

$ws.Range("A3").Value = "2025-09-06 12:31:29"
$ws.Range("A4").Value = "2025-09-06 12:31:29"
$ws.Range("A5").Value = "2025-09-06 12:31:29"
$ws.Range("A6").Value = "2025-09-06 12:31:29"

# Column H ("スキル概要") is now wider to fit the longer tag text.
# (ColumnWidth as read/written by this host is offset by 5/6 of a character
# from the stored OOXML <col width> value, so subtract that to land on 21.)
$ws.Columns.Item(8).ColumnWidth = 21 - 5/6

# Rebuild the hyperlinks collection so the rId -> URL mapping follows the
# rows that shifted down, plus the brand-new link for the inserted row.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5388329")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5388066")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5388189")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5385681")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5388228")
